$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New model-name order for column A (rows 2-26), a permutation of the
# existing model_6_7_* identifiers.
$names = @(
    "model_6_7_0",
    "model_6_7_22",
    "model_6_7_21",
    "model_6_7_20",
    "model_6_7_19",
    "model_6_7_18",
    "model_6_7_17",
    "model_6_7_16",
    "model_6_7_15",
    "model_6_7_14",
    "model_6_7_13",
    "model_6_7_23",
    "model_6_7_12",
    "model_6_7_10",
    "model_6_7_9",
    "model_6_7_8",
    "model_6_7_7",
    "model_6_7_6",
    "model_6_7_5",
    "model_6_7_4",
    "model_6_7_3",
    "model_6_7_2",
    "model_6_7_1",
    "model_6_7_11",
    "model_6_7_24"
)

# New metric values (same for every row, columns B..Q).
$values = @(
    [double]"0.9999632043074972",
    [double]"0.9992467334840777",
    [double]"0.999998685648971",
    [double]"0.9999999999999448",
    [double]"0.999999527764742",
    [double]"3.434716380714357e-05",
    [double]"0.0007031412280348575",
    [double]"9.346684710893685e-07",
    [double]"6.78732299646294e-14",
    [double]"4.673342694812993e-07",
    [double]"0.0002886930549261025",
    [double]"0.005860645340501639",
    [double]"1.000883096620067",
    [double]"0.006110145145855123",
    [double]"70.55798222036485",
    [double]"101.0298778420699"
)

for ($i = 0; $i -lt 25; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    for ($c = 0; $c -lt $values.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $values[$c]
    }
}

$wb.Save()
